$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.450.22"
$ws.Range("E2").Value = "  +0.17%  "

# Row 3
$ws.Range("D3").Value = "1.947.57"
$ws.Range("E3").Value = "  -1.65%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.51"
$ws.Range("E5").Value = "  -0.53%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("E6").Value = "  -1.94%  "

# Row 7
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.88"
$ws.Range("E8").Value = "  -2.59%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.363"
$ws.Range("E9").Value = "  -3.56%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0847"
$ws.Range("E10").Value = "  +4.12%  "

# Row 11
$ws.Range("E11").Value = "  +0.15%  "

# Row 12
$ws.Range("D12").Value = "2.234.38"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.815"
$ws.Range("E13").Value = "  -5.20%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.16"
$ws.Range("E14").Value = "  -9.96%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.51"
$ws.Range("E15").Value = "  -3.25%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.17"
$ws.Range("E16").Value = "  -4.88%  "

# Row 17
$ws.Range("D17").Value = "1.928.85"
$ws.Range("E17").Value = "  -2.39%  "

# Row 18
$ws.Range("D18").Value = "36.337.02"
$ws.Range("E18").Value = "  -0.05%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0875"
$ws.Range("E19").Value = "  +1.60%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.53"
$ws.Range("E20").Value = "  -1.40%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "229.22"
$ws.Range("E21").Value = "  -2.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.02"
$ws.Range("E22").Value = "  -5.30%  "

# Row 23
$ws.Range("E23").Value = "  -0.13%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  -7.81%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.28"
$ws.Range("E25").Value = "  -0.75%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.18"
$ws.Range("E26").Value = "  -8.68%  "

# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.13"
$ws.Range("E27").Value = "  -0.26%  "

# Row 28
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.137"
$ws.Range("E28").Value = "  +8.16%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.26"
$ws.Range("E29").Value = "  -2.54%  "

# Row 30
$ws.Range("E30").Value = "  -1.75%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.14"
$ws.Range("E31").Value = "  -4.58%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.61"
$ws.Range("E32").Value = "  -5.54%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0638"
$ws.Range("E33").Value = "  +1.33%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.25"
$ws.Range("E34").Value = "  -3.38%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.16"
$ws.Range("E35").Value = "  -1.87%  "

# Row 36
$ws.Range("E36").Value = "  +0.01%  "

# Row 37
$ws.Range("E37").Value = "  +1.37%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.13"
$ws.Range("E38").Value = "  -5.79%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.02"
$ws.Range("E39").Value = "  -1.71%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0972"
$ws.Range("E40").Value = "  +1.05%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.89"
$ws.Range("E41").Value = "  -0.22%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0211"
$ws.Range("E42").Value = "  -1.06%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.16"
$ws.Range("E43").Value = "  -6.82%  "

# Row 44
$ws.Range("D44").Value = "1.354.56"
$ws.Range("E44").Value = "  -1.23%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.57"
$ws.Range("E45").Value = "  -4.04%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.02"
$ws.Range("E46").Value = "  -6.45%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.14"
$ws.Range("E47").Value = "  -5.55%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.09"
$ws.Range("E48").Value = "  -6.08%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.84"
$ws.Range("E49").Value = "  -0.54%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.41"
$ws.Range("E50").Value = "  -1.32%  "

# Row 51
$ws.Range("D51").Value = "2.125.48"
$ws.Range("E51").Value = "  -1.95%  "

